$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- New journal entries (rows 52-53) ---
$ws.Range("A52").Value = 44697
$ws.Range("B52").Value = 0.40972222222222227
$ws.Range("C52").Value = 0.4375
$ws.Range("E52").Value = "Fin de la doc concernant `nle placement de 1 régiment"

$ws.Range("A53").Value = 44697
$ws.Range("B53").Value = 0.4375
$ws.Range("C53").Value = 0.51041666666666663
$ws.Range("E53").Value = "Algorithme avec plusieurs régiment(Placement)"

# --- View state: scroll down to row 49 and select E53 ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E53").Select()
